# Generate Report for Handback
# The handback transform failed for the file
# 7b1e68da-4e50-4840-9a71-7a6a9cb3827a (row 3 of the zh-cn / de-de / Overview
# tables), so its status is updated and an error detail message is recorded
# for both locales.

$wb = $excel.ActiveWorkbook

$statusFailed = "Handback transform failed"

# --- Overview sheet: update the zh-cn / de-de status columns for row 3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusFailed
$wsOverview.Range("F3").Value = $statusFailed

# --- zh-cn sheet: update Status + Error Detail for row 3, widen column P ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusFailed
$wsZhCn.Range("P3").Value = "Handback file name: ombvs31p.tgs is different with handoff file name: 7b1e68da-4e50-4840-9a71-7a6a9cb3827a.5236d049c62bd9595833ccbf10469b21ddf0aff8.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: update Status + Error Detail for row 3, widen column P ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusFailed
$wsDeDe.Range("P3").Value = "Handback file name: ombvs31p.tgs is different with handoff file name: 7b1e68da-4e50-4840-9a71-7a6a9cb3827a.5236d049c62bd9595833ccbf10469b21ddf0aff8.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
